# Applies the "escape_tpl" template update:
#  - the "{{r page_break }} " paragraph loses its stray trailing space
#  - a new "{{ new_listing }}" paragraph (italic, accent6 theme color,
#    dark-yellow highlight) is inserted right after it, before "END"
#  - the section's page orientation is stamped explicitly (portrait)

$d = $word.ActiveDocument

# --- locate the "{{r page_break }} " paragraph -----------------------
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text -like "*page_break*") {
        $target = $cand
    }
}

$r = $target.Range

# Drop the single trailing space before the paragraph mark:
# "{{r page_break }} " -> "{{r page_break }}"
$endOfPara = $r.End
$trailing = $d.Range($endOfPara - 2, $endOfPara - 1)
if ($trailing.Text -eq " ") {
    $trailing.Text = ""
}

# --- insert the new "{{ new_listing }}" paragraph right after -------
$afterPara = $target.Range
$newPara = $afterPara.InsertParagraphAfter()

$newParaObj = $target.Next()
$newRange = $newParaObj.Range
$newRange.InsertAfter("{{ new_listing }}")

$newRange.Font.Italic = $true
$newRange.Font.ItalicBi = $true
$newRange.Font.TextColor.ObjectThemeColor = 9
$newRange.HighlightColorIndex = 14

# --- explicitly stamp the section's orientation (portrait) ----------
$d.Sections.Item(1).PageSetup.Orientation = 0
